# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet (fund holdings detail) right before
#    the existing "总计" (totals) sheet -- built by duplicating "总计" (so
#    it inherits the same header/border styling) and then overwriting its
#    contents.
# 2. Insert a new top data row into "总计" summarizing the 2022-Q1 quarter
#    (date/count/value), pushing the existing 2021-Q4 / 2021-Q3 rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell as TEXT (no numeric auto-coercion,
# no left-over custom NumberFormat/style) -- matches how the other
# "XXXX-QX" detail sheets store their numeric-looking text cells
# (e.g. "012262", "16.19", "0.2590").
# ---------------------------------------------------------------------
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Duplicate "总计" to inherit its header/border formatting, place the
#    duplicate right before it, rename it to "2022-Q1" and replace the
#    data with the new quarter's fund holdings.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$totalIndex = $total.Index
$total.Copy($total)

# NOTE: this COM host resolves worksheet object references by their
# current tab position, not stable identity. After Copy() inserts the
# duplicate *before* the original, the position that $total used to
# occupy now holds the new duplicate, and the original "总计" sheet has
# been pushed one slot to the right -- so both sheets must be re-fetched
# by position/name rather than reusing the old $total handle.
$q1 = $wb.Worksheets.Item($totalIndex)
$total = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Clear the old (copied) rows before writing the new layout.
$q1.Range("A1:D3").ClearContents()

# The copied "总计" sheet only had styled cells out to column D / row 3;
# extend the same header (B1:D1 style) and row-number (A2 style) look to
# the extra columns/rows this sheet needs.
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)
$q1.Range("A2").Copy()
$q1.Range("A3:A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data rows
$q1Data = @(
    @("012262", "华宝可持续发展混合A", "16.19", "66.51", "1.60", "0.2590", 8),
    @("012263", "华宝可持续发展混合C", "6.42",  "66.51", "1.60", "0.1027", 8),
    @("011734", "国寿安保裕丰混合型证券投资基金A", "5.01", "20.33", "0.54", "0.0271", 7),
    @("011735", "国寿安保裕丰混合型证券投资基金C", "1.37", "20.33", "0.54", "0.0074", 7)
)

$r = 2
foreach ($row in $q1Data) {
    $q1.Range("A$r").Value = $($r - 2)
    Set-TextValue $q1.Range("B$r") $row[0]
    $q1.Range("C$r").Value = $row[1]
    Set-TextValue $q1.Range("D$r") $row[2]
    Set-TextValue $q1.Range("E$r") $row[3]
    Set-TextValue $q1.Range("F$r") $row[4]
    Set-TextValue $q1.Range("G$r") $row[5]
    $q1.Range("H$r").Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Insert a new summary row for 2022-Q1 at the top of "总计"'s data
#    (row 2, just below the header), shifting the rest down.
# ---------------------------------------------------------------------
$total.Range("A2:D2").Insert()

# Insert() leaves the new B2:D2 cells carrying a blended style picked up
# from the row above (the bold header) -- reset to plain/unstyled like
# the other data cells, and give A2 the same "row index" style (s=2) as
# A3/A4 by copying it across.
$total.Range("B2:D2").Style = "Normal"
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.4

# Re-sequence the row-index column (A) -- Insert() shifts existing cells
# down verbatim, it doesn't renumber them.
$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
